$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E: "Nullable Enums" - add header with the same style as the
# other header cells (copy formatting from D1, then set the text).
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Nullable Enums"

# Fill in the nullable-enum values. Row 5 is intentionally left blank
# (null) to exercise the "nullable" part of the new column.
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 1

# Match the new selection left behind in the saved workbook.
[void]$ws.Range("E9").Select()
